# Convert M2Doc "Word field" style tokens ({ m:'An empty table'.emptyTable() }
# stored as a real Word field with fldChar begin/instrText/fldChar end) into
# the plain-text curly-brace notation used by the
# TokenIteratorFieldRewriterSplit parser: four separate <w:t> runs holding
# "{", "m", the rest of the expression, and "}".

$d = $word.ActiveDocument

# Walk the Fields collection back-to-front so deleting/inserting around one
# field never invalidates the indices of the fields we haven't visited yet.
for ($i = $d.Fields.Count; $i -ge 1; $i--) {

    $fld = $d.Fields.Item($i)
    $code = $fld.Code.Text

    # Only touch M2Doc fields, i.e. field codes that look like
    # " m:....... " (a literal leading/trailing space around the expression).
    if ($code.Length -lt 3 -or $code.Substring(0, 1) -ne " " -or $code.Substring($code.Length - 1, 1) -ne " ") {
        continue
    }

    # Strip the one leading/trailing space that the field-code convention
    # adds (it plays the same role the '{' and '}' braces play in the
    # plain-text notation), then split "m" (the binding marker) from the
    # rest of the expression - matching the original run boundaries.
    $inner = $code.Substring(1, $code.Length - 2)
    $mPart = $inner.Substring(0, 1)
    $restPart = $inner.Substring(1, $inner.Length - 1)

    # Remember where the field starts so we can re-insert plain text at the
    # same spot, then remove the field (begin/instrText*/end) entirely.
    $insertPos = $fld.Code.Start - 1
    $fld.Delete()

    # Build a minimal WordprocessingML package fragment with exactly the
    # four runs required, and drop it in via InsertXML so the run
    # boundaries are preserved verbatim (plain Range.InsertAfter calls get
    # coalesced into a single run on save).
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        '<w:p>' + `
            '<w:r><w:t>{</w:t></w:r>' + `
            '<w:r><w:t>' + $mPart + '</w:t></w:r>' + `
            '<w:r><w:t>' + $restPart + '</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r = $d.Range($insertPos, $insertPos)
    $r.InsertXML($xml)
}
